$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44326, 0, 5, 94.6969696969697),
    @(44327, 1, 6, 113.6363636363636),
    @(44328, 0, 6, 113.6363636363636),
    @(44329, 0, 4, 75.75757575757575)
)

$row = 252
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Copy style from A251 to the new date cells A252:A255
$ws.Range("A251").Copy()
$ws.Range("A252:A255").PasteSpecial(-4122)  # xlPasteFormats
